$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old sub-header row (Hiver/Eté/Année + units), shifting data rows up by one.
$ws.Rows(2).Delete()

# Rebuild the header row with the new "idx / idx2 / Name / Date Start / Date End" columns
# plus the renamed measurement columns.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
# E1 previously held a leftover styled header ("mation") - clear its style so the
# new "Date End" header lands on the default (unstyled) xf, matching A1:D1.
$ws.Range("E1").ClearFormats()
$ws.Range("E1").Value = "Date End"

$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# Give the measurement headers (F1:K1) a dedicated style (Arial 9, general format),
# mirroring the existing data-cell font without inheriting the data cells' number
# formats. A scratch named style is used to synthesize the xf and then discarded so
# the workbook keeps only the default three built-in cell styles.
$tempStyle = $wb.Styles.Add("TempHeaderStyle")
$tempStyle.Font.Name = "Arial"
$tempStyle.Font.Size = 9

$measureHeaders = $ws.Range("F1:K1")
$measureHeaders.Style = "TempHeaderStyle"

$tempStyle.Delete() | Out-Null

$ws.Range("A2:K2").Select() | Out-Null
